# Rename the Korean activity names in column A to the English
# placeholders a..j, move the active selection to A12, and configure
# the page setup (A4 paper, portrait orientation) on the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("activity_list")

$names = @("a", "b", "c", "d", "e", "f", "g", "h", "i", "j")
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $names[$i]
}

$ws.Range("A12").Select()

$ws.PageSetup.PaperSize = 9   # xlPaperA4
$ws.PageSetup.Orientation = 1 # xlPortrait
